# Adapt column header formatting to respective input file names (#7)
# - rename the "*_old" headers to "*_FV2210"
# - rename the "*_new" headers to "*_FV2304"
# - turn the data range into an Excel Table (Table1)
# - freeze the header row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Rename header row (row 1) -----------------------------------------
$oldSuffixCols = @("A","B","C","D","E","F","G","H","I","J")
$newSuffixCols = @("L","M","N","O","P","Q","R","S","T","U")
$baseNames = @("Segmentname","Segmentgruppe","Segment","Datenelement","Segment ID","Code","Qualifier","Beschreibung","Bedingungsausdruck","Bedingung")

for ($i = 0; $i -lt $baseNames.Length; $i++) {
    $ws.Range($oldSuffixCols[$i] + "1").Value = $baseNames[$i] + "_FV2210"
    $ws.Range($newSuffixCols[$i] + "1").Value = $baseNames[$i] + "_FV2304"
}
# column K ("diff") is unchanged

# --- 2) Turn the used range into a table ------------------------------------
$dataRange = $ws.Range("A1:U70")
$tbl = $ws.ListObjects.Add(1, $dataRange, 0, 1)
$tbl.Name = "Table1"

# --- 3) Freeze the header row ------------------------------------------------
$ws.Range("A2").Select()
$win = $excel.ActiveWindow
$win.FreezePanes = $true
